$d = $word.ActiveDocument

# Helper: replace the text of the Nth paragraph (by 1-based occurrence count)
# whose trimmed text exactly equals $oldText, with $newText. Using a direct
# Range.Text assignment (rather than Find.Execute's Replace) avoids invoking
# Word's "replace straight quotes with smart quotes" AutoCorrect/AutoFormat
# behavior, which only triggers on simulated typing/Find-Replace operations.
function Replace-ParaText($oldText, $newText, $occurrence) {
    $seen = 0
    foreach ($p in $d.Paragraphs) {
        $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $oldText) {
            $seen = $seen + 1
            if ($seen -eq $occurrence) {
                $r = $p.Range
                $r.MoveEnd(1, -1) | Out-Null
                $r.Text = $newText
                return $true
            }
        }
    }
    return $false
}

Replace-ParaText "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations" "• Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions" 1 | Out-Null
Replace-ParaText "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics" "• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics and demographic analysis" 1 | Out-Null
Replace-ParaText "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" "• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" 1 | Out-Null
Replace-ParaText "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" "• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" 1 | Out-Null
Replace-ParaText "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" "• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" 1 | Out-Null
Replace-ParaText "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" "• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" 1 | Out-Null
Replace-ParaText "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" "• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" 1 | Out-Null
Replace-ParaText "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" "• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" 1 | Out-Null
Replace-ParaText "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" "• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" 1 | Out-Null
Replace-ParaText "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products" "• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products" 1 | Out-Null
Replace-ParaText "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" "• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" 1 | Out-Null
Replace-ParaText "• Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies" "• Architected and developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies" 1 | Out-Null
Replace-ParaText "• Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu" "• Collaborated with senior officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to enhance platform using Docker and Ubuntu" 1 | Out-Null
Replace-ParaText "• Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions" "• Conceived and developed predictive application using Python, Pandas, and Jupyter to forecast how crisis economies respond to different humanitarian interventions" 1 | Out-Null
Replace-ParaText "• Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting" "• Architected and developed RACSO, a comprehensive web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting" 1 | Out-Null
Replace-ParaText "• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner" "• Led RFP process and analyzed bids from 1,200 vendors before selecting optimal implementation partner" 1 | Out-Null
Replace-ParaText "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research" "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions" 1 | Out-Null
Replace-ParaText "• Designed survey deployment system facilitating thousands of simultaneous phone surveys" "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs" 1 | Out-Null
Replace-ParaText "• Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system" "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system" 1 | Out-Null
Replace-ParaText "• Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill" "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill" 1 | Out-Null
Replace-ParaText "• Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization" "• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO" 1 | Out-Null
Replace-ParaText "• Made all technology decisions and practices for massive multinational non-governmental organization" "• Directed all technology decisions and practices for massive multinational non-governmental organization" 1 | Out-Null
Replace-ParaText "• Wrote comprehensive frameworks for internal and external technology audits" "• Developed comprehensive frameworks for internal and external technology audits" 1 | Out-Null
Replace-ParaText "• Trained beneficiaries on spatial and Census data analysis for public health research" "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research" 1 | Out-Null
Replace-ParaText "• Trained NGO staff in web development using Drupal, PHP, and MySQL" "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL" 1 | Out-Null
Replace-ParaText "Political Research and Data Analysis" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns" 2 | Out-Null  # 2nd occurrence only; 1st (RESEARCH DIRECTOR section) stays as-is
Replace-ParaText "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions" 1 | Out-Null
Replace-ParaText "Political Field Operations and Data Management" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns" 1 | Out-Null
Replace-ParaText "• Administered all quantitative and qualitative research operations ensuring reporting accuracy" "• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions" 1 | Out-Null
Replace-ParaText "• Managed comprehensive survey fielding for multi-million dollar research firm" "• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm" 1 | Out-Null
